$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.851.64"
$ws.Range("E2").Value = "  +6.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.92"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.70"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.41"
$ws.Range("E6").Value = "  +13.51%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  +8.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.95"
$ws.Range("E10").Value = "  +13.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0803"
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").Value = "  +7.13%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.659.51"
$ws.Range("E14").Value = "  +3.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.305.86"
$ws.Range("E15").Value = "  +4.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.02"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("E17").Value = "  +6.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.843.17"
$ws.Range("E18").Value = "  +6.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("E19").Value = "  +25.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("E20").Value = "  +4.64%  "
$ws.Range("E21").Value = "  +4.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.89"
$ws.Range("E22").Value = "  +3.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.61"
$ws.Range("E23").Value = "  +5.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  +5.72%  "
$ws.Range("E25").Value = "  +6.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "44.22"
$ws.Range("E27").Value = "  +16.60%  "
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.96"
$ws.Range("E29").Value = "  +6.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.23"
$ws.Range("E30").Value = "  +4.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.80"
$ws.Range("E31").Value = "  +8.19%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("E32").Value = "  +9.05%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0810"
$ws.Range("E33").Value = "  +8.29%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.88"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("E35").Value = "  +14.91%  "
$ws.Range("E36").Value = "  +12.41%  "
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("E38").Value = "  +8.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.21"
$ws.Range("E39").Value = "  +24.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.03"
$ws.Range("E40").Value = "  +15.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  +10.18%  "
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("E43").Value = "  +14.92%  "
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.861.77"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.86"
$ws.Range("E46").Value = "  +20.27%  "
$ws.Range("E47").Value = "  +10.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "74.87"
$ws.Range("E48").Value = "  +13.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.91"
$ws.Range("E49").Value = "  +11.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.36"
$ws.Range("E50").Value = "  +3.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.06"
$ws.Range("E51").Value = "  +5.81%  "
